$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 ("ჩართულ ბავშვთა (აღსაზრდელები) რაოდენობა (კაცი)") held stray
# percentage-like decimals that were mistakenly copy/pasted from the
# "average children per teacher" row. Replace them with the correct
# integer headcounts and format the row as whole numbers, matching the
# rest of the sheet's integer rows.
$ws.Range("B5").Value = 661
$ws.Range("C5").Value = 619
$ws.Range("D5").Value = 650
$ws.Range("E5").Value = 657
$ws.Range("F5").Value = 616

$ws.Range("B5:F5").NumberFormat = "#\ ##0"

# Restore the active selection to where the edit was made.
$ws.Range("B5").Select()
